# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape. Cell values that look numeric (e.g. "1.000", "0.06790")
# are written with a leading apostrophe so Excel keeps them as literal
# text (matching the workbook's inlineStr storage) instead of silently
# coercing them into numbers and dropping significant trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.728.83'
$ws.Range("E2").Value = '  +2.27%  '
$ws.Range("D3").Value = '1.891.81'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = "'244.95"
$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = "'0.4917"
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = "'0.2959"
$ws.Range("E8").Value = '  +0.98%  '
$ws.Range("D9").Value = "'0.06790"
$ws.Range("E9").Value = '  +2.88%  '
$ws.Range("D10").Value = '1.889.08'
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("E11").Value = '  +3.97%  '
$ws.Range("D12").Value = "'0.07229"
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("D13").Value = "'91.32"
$ws.Range("E13").Value = '  +5.89%  '
$ws.Range("D14").Value = "'0.6794"
$ws.Range("E14").Value = '  +1.83%  '
$ws.Range("D15").Value = "'5.049"
$ws.Range("E15").Value = '  +2.71%  '
$ws.Range("D16").Value = '30.683.88'
$ws.Range("E16").Value = '  +2.28%  '
$ws.Range("D17").Value = "'0.000008003"
$ws.Range("E17").Value = '  +2.69%  '
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = '  +3.11%  '
$ws.Range("D20").Value = '2.131.43'
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("D23").Value = "'191.72"
$ws.Range("E23").Value = '  +33.65%  '
$ws.Range("D24").Value = "'6.121"
$ws.Range("E24").Value = '  +4.50%  '
$ws.Range("D25").Value = "'9.374"
$ws.Range("E25").Value = '  +3.05%  '
$ws.Range("D26").Value = "'154.80"
$ws.Range("E26").Value = '  +2.28%  '
$ws.Range("D27").Value = "'19.18"
$ws.Range("E27").Value = '  +13.26%  '
$ws.Range("D28").Value = "'1.906"
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("D29").Value = "'1.401"
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("D30").Value = "'4.346"
$ws.Range("E30").Value = '  +3.68%  '
$ws.Range("D31").Value = "'0.09097"
$ws.Range("E31").Value = '  +4.01%  '
$ws.Range("D32").Value = "'4.019"
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("D33").Value = "'0.05206"
$ws.Range("E33").Value = '  +3.74%  '
$ws.Range("D34").Value = "'0.7538"
$ws.Range("E34").Value = '  +5.29%  '
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("D36").Value = "'2.771"
$ws.Range("E36").Value = '  +4.13%  '
$ws.Range("D37").Value = "'0.01847"
$ws.Range("E37").Value = '  +1.61%  '
$ws.Range("D38").Value = "'2.687"
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").Value = "'2.149"
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("D40").Value = "'0.9377"
$ws.Range("E40").Value = '  +0.71%  '
$ws.Range("D41").Value = "'0.4427"
$ws.Range("E41").Value = '  +4.96%  '
$ws.Range("D42").Value = "'105.27"
$ws.Range("E42").Value = '  +2.10%  '
$ws.Range("D43").Value = "'5.771"
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("D44").Value = "'0.9999"
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").Value = "'7.615"
$ws.Range("E45").Value = '  +3.36%  '
$ws.Range("D46").Value = "'0.1348"
$ws.Range("E46").Value = '  +6.33%  '
$ws.Range("D47").Value = "'0.05865"
$ws.Range("E47").Value = '  +2.92%  '
$ws.Range("D48").Value = "'8.732"
$ws.Range("E48").Value = '  +5.97%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = "'1.426"
$ws.Range("E49").Value = '  +6.57%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").Value = "'0.3944"
$ws.Range("E50").Value = '  +4.76%  '
$ws.Range("D51").Value = "'33.66"
$ws.Range("E51").Value = '  +2.81%  '
